# Insert a new data row at row 156 (pushing existing rows 156-248 down to 157-249)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 156; this shifts row 156..248 down to 157..249
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record.
$ws.Cells.Item(156, 1).Value = 10
$ws.Cells.Item(156, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(156, 3).Value = "La Araucanía"
$ws.Cells.Item(156, 4).Value = 44572
$ws.Cells.Item(156, 5).Value = 9
$ws.Cells.Item(156, 6).Value = 100112044
$ws.Cells.Item(156, 7).Value = "Perejil"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 65
$ws.Cells.Item(156, 11).Value = 4000
$ws.Cells.Item(156, 12).Value = 4000
$ws.Cells.Item(156, 13).Value = 4000
$ws.Cells.Item(156, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(156, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(156, 16).Value = 1333
$ws.Cells.Item(156, 17).Value = 3
$ws.Cells.Item(156, 18).Value = "Hortaliza"
